# Adicionado Scaller para garantir que valores fiquem entre min e max.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Discreto 2" (row 5): Max 5 -> 10, Min 3 -> 2, Step 0.5 -> 1
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1

# New variable "Temperatura" added as row 6
$ws.Range("A6").Value = "Temperatura"
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = 50
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = 0.95
$ws.Range("H6").Value = "Continuous"
